$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "QF13002"
$ws.Range("B3").Value = "Aristides"
$ws.Range("C3").Value = "Fuentes"
$ws.Range("D3").Value = 28
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 7.5
$ws.Range("G3").Value = 2013
$ws.Range("H3").Value = 7
$ws.Range("I3").Value = 1

$ws.Range("I3").Select() | Out-Null
